$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

$ws.Range("C2").Value = "firefox"
$ws.Range("C4").Value = "firefox"

$ws.Range("C4").Select()
